$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "A2"  = "to open (something) (～を)"
    "A3"  = "to close (something) (～を)"
    "A4"  = "to teach; to instruct (person に thing を)"
    "A5"  = "to forget; to leave behind (～を)"
    "A6"  = "to get off (～を)"
    "A7"  = "to borrow (person に thing を)"
    "A9"  = "to turn on (～を)"
    "A10" = "to call (～に)"
    "A11" = "to bring (a person) (～を)"
    "A12" = "to bring (a thing) (～を)"
    "A46" = "to return (a thing) (person に thing を)"
    "A47" = "to turn off; to erase (～を)"
    "A49" = "to sit down (seat に)"
    "A52" = "to use (～を)"
    "A53" = "to help (person/task を)"
    "A54" = "to enter (～に)"
    "A55" = "to carry; to hold (～を)"
    "A56" = "(1) to be absent (from...) (～を); (2) to rest"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
